# Scheduled runner update: refresh market-board derived price/profit columns
# (H:N) across the Leve profit sheets with the latest fetched values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1972.3077
$ws.Range("I38").Value = 55
$ws.Range("J38").Value = 5040
$ws.Range("K38").Value = 165
$ws.Range("L38").Value = 15120
$ws.Range("M38").Value = 207
$ws.Range("N38").Value = -15864
$ws.Range("H87").Value = 67499.5
$ws.Range("J87").Value = 99999
$ws.Range("L87").Value = 99999
$ws.Range("N87").Value = -102495
$ws.Range("H90").Value = 67499.5
$ws.Range("J90").Value = 99999
$ws.Range("L90").Value = 299997
$ws.Range("N90").Value = -312477
$ws.Range("H103").Value = 508.42105
$ws.Range("I103").Value = 447.66666
$ws.Range("J103").Value = 736.25
$ws.Range("K103").Value = 1342.99998
$ws.Range("L103").Value = 2208.75
$ws.Range("M103").Value = -756.9999800000001
$ws.Range("N103").Value = -3380.75
$ws.Range("H112").Value = 1393.9736
$ws.Range("J112").Value = 1280.3438
$ws.Range("L112").Value = 3841.0314
$ws.Range("N112").Value = -6057.0314
$ws.Range("H127").Value = 670.5
$ws.Range("I127").Value = 629.3333
$ws.Range("K127").Value = 1887.9999
$ws.Range("M127").Value = 3072.0001
$ws.Range("H132").Value = 1699.8983
$ws.Range("I132").Value = 1077.3396
$ws.Range("K132").Value = 3232.0188
$ws.Range("M132").Value = -702.0187999999998
$ws.Range("H135").Value = 733.6667
$ws.Range("I135").Value = 733.6667
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6603.0003
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -4068.0003
$ws.Range("H137").Value = 38463348
$ws.Range("I137").Value = 55556610
$ws.Range("K137").Value = 166669830
$ws.Range("M137").Value = -166667280
$ws.Range("H138").Value = 4011.6316
$ws.Range("J138").Value = 5291.758
$ws.Range("L138").Value = 15875.274
$ws.Range("N138").Value = -26155.274

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 22630
$ws.Range("J76").Value = 34999
$ws.Range("L76").Value = 34999
$ws.Range("N76").Value = -35675
$ws.Range("H79").Value = 22630
$ws.Range("J79").Value = 34999
$ws.Range("L79").Value = 34999
$ws.Range("N79").Value = -37339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 27519.092
$ws.Range("I99").Value = 11587.375
$ws.Range("J99").Value = 70003.664
$ws.Range("K99").Value = 11587.375
$ws.Range("L99").Value = 70003.664
$ws.Range("M99").Value = -10089.375
$ws.Range("N99").Value = -72999.664
$ws.Range("H140").Value = 99998.5
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6619.5625
$ws.Range("I31").Value = 5429.28
$ws.Range("K31").Value = 5429.28
$ws.Range("M31").Value = -5134.28
$ws.Range("H34").Value = 6619.5625
$ws.Range("I34").Value = 5429.28
$ws.Range("K34").Value = 5429.28
$ws.Range("M34").Value = -5227.28
$ws.Range("H134").Value = 30364260
$ws.Range("I134").Value = 33400286
$ws.Range("K134").Value = 100200858
$ws.Range("M134").Value = -100198323

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 7261.3335
$ws.Range("J82").Value = 9492
$ws.Range("L82").Value = 28476
$ws.Range("N82").Value = -29288
$ws.Range("H85").Value = 7261.3335
$ws.Range("J85").Value = 9492
$ws.Range("L85").Value = 28476
$ws.Range("N85").Value = -31284
$ws.Range("H93").Value = 5511.1665
$ws.Range("J93").Value = 5830.4546
$ws.Range("L93").Value = 17491.3638
$ws.Range("N93").Value = -21235.3638
$ws.Range("H122").Value = 290.2353
$ws.Range("I122").Value = 357
$ws.Range("J122").Value = 253.81818
$ws.Range("K122").Value = 3213
$ws.Range("L122").Value = 2284.36362
$ws.Range("M122").Value = -763
$ws.Range("N122").Value = -7184.36362
$ws.Range("H131").Value = 1633.625
$ws.Range("I131").Value = 966.6667
$ws.Range("J131").Value = 1687.7028
$ws.Range("K131").Value = 2900.0001
$ws.Range("L131").Value = 5063.1084
$ws.Range("M131").Value = 2139.9999
$ws.Range("N131").Value = -15143.1084

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3991.25
$ws.Range("I126").Value = 5532
$ws.Range("J126").Value = 3290.9092
$ws.Range("K126").Value = 16596
$ws.Range("L126").Value = 9872.7276
$ws.Range("M126").Value = -14126
$ws.Range("N126").Value = -14812.7276
$ws.Range("H132").Value = 27860.195
$ws.Range("I132").Value = 28133.768
$ws.Range("K132").Value = 84401.304
$ws.Range("M132").Value = -81871.304

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1009.5789
$ws.Range("I22").Value = 1002.125
$ws.Range("J22").Value = 1049.3334
$ws.Range("K22").Value = 1002.125
$ws.Range("L22").Value = 1049.3334
$ws.Range("M22").Value = -707.125
$ws.Range("N22").Value = -1639.3334
$ws.Range("H27").Value = 1009.5789
$ws.Range("I27").Value = 1002.125
$ws.Range("J27").Value = 1049.3334
$ws.Range("K27").Value = 1002.125
$ws.Range("L27").Value = 1049.3334
$ws.Range("M27").Value = -895.125
$ws.Range("N27").Value = -1263.3334
$ws.Range("H46").Value = 5244.278
$ws.Range("I46").Value = 8640.556
$ws.Range("K46").Value = 8640.556
$ws.Range("M46").Value = -8452.556
$ws.Range("H100").Value = 7273.364
$ws.Range("I100").Value = 2854.1765
$ws.Range("J100").Value = 22298.6
$ws.Range("K100").Value = 2854.1765
$ws.Range("L100").Value = 22298.6
$ws.Range("M100").Value = -2313.1765
$ws.Range("N100").Value = -23380.6
$ws.Range("H132").Value = 2608.8438
$ws.Range("I132").Value = 2579.7856
$ws.Range("J132").Value = 2812.25
$ws.Range("K132").Value = 7739.3568
$ws.Range("L132").Value = 8436.75
$ws.Range("M132").Value = -5209.3568
$ws.Range("N132").Value = -13496.75
$ws.Range("H136").Value = 3964.4
$ws.Range("I136").Value = 2455.5
$ws.Range("K136").Value = 7366.5
$ws.Range("M136").Value = -4816.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9528297
$ws.Range("I132").Value = 11768485
$ws.Range("J132").Value = 7499.25
$ws.Range("K132").Value = 35305455
$ws.Range("L132").Value = 22497.75
$ws.Range("M132").Value = -35302925
$ws.Range("N132").Value = -27557.75
